# Update the two-digit division worksheet: replace each problem's
# "a÷b=" text with the new problem text, cell by cell, in document order.
#
# Note: "19÷5=" occurs twice in the source document with two different
# replacements ("21÷8=" then "50÷2="). Using MatchWholeWord Find/Replace
# with Replace=wdReplaceOne (1) on $d.Content for each pair - in the same
# order the pairs appear in the document - always targets the first
# remaining (left-to-right, top-to-bottom) occurrence, so the duplicate
# resolves correctly without disturbing any other cell.

$d = $word.ActiveDocument

$replacements = @(
    @("33÷8=", "26÷8="),
    @("83÷5=", "30÷3="),
    @("74÷4=", "17÷6="),
    @("48÷6=", "75÷5="),
    @("79÷2=", "27÷6="),
    @("53÷2=", "89÷9="),
    @("37÷3=", "15÷3="),
    @("70÷3=", "50÷4="),
    @("67÷9=", "72÷8="),
    @("29÷8=", "36÷6="),
    @("29÷4=", "36÷9="),
    @("81÷7=", "62÷3="),
    @("26÷4=", "65÷8="),
    @("67÷7=", "54÷8="),
    @("19÷5=", "21÷8="),
    @("25÷2=", "76÷5="),
    @("55÷2=", "83÷4="),
    @("93÷5=", "63÷6="),
    @("93÷2=", "16÷8="),
    @("25÷8=", "69÷5="),
    @("89÷5=", "66÷5="),
    @("19÷5=", "50÷2="),
    @("95÷5=", "61÷7="),
    @("61÷4=", "36÷3="),
    @("11÷9=", "19÷6=")
)

$count = 0
foreach ($pair in $replacements) {
    $oldText = $pair[0]
    $newText = $pair[1]
    $found = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 1)
    if ($found) {
        $count = $count + 1
    }
}
Write-Output "Replacements applied: $count of $($replacements.Count)"
